$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 21 de Marzo de 2020 a las 05:16"

$ws.Range("A8").Value = "Estados Unidos"
$ws.Range("A9").Value = "Iran"

$ws.Range("B8").Value = 19652
$ws.Range("C8").Value = 269
$ws.Range("D8").Value = 147
$ws.Range("E8").Value = 19241
$ws.Range("F8").Value = 64
$ws.Range("G8").Value = 8
$ws.Range("H8").Value = 264

$ws.Range("B9").Value = 19644
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 6745
$ws.Range("E9").Value = 11466
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 0
$ws.Range("H9").Value = 1433

$ws.Range("E40").Value = 252
$ws.Range("G40").Value = 2
$ws.Range("H40").Value = 2
